$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new columns: G "max clique", I "odd cycle" -------------------------
$ws.Range("G1").Value = "max clique"
$ws.Range("I1").Value = "odd cycle"

# Column G width so the exported OOXML width matches the original (10 chars)
$ws.Columns("G").ColumnWidth = 9.140625

# --- data for rows 23..42 ------------------------------------------------
# max-clique size found for each graph
$G = @{
    23 = 3
    24 = 3
    25 = 4
    26 = 4
    27 = 2
    28 = 3
    29 = 8
    30 = 6
    31 = 3
    32 = 3
    33 = 5
    34 = 2
    35 = 7
    36 = 3
    37 = 5
    38 = 3
    39 = 6
    40 = 6
    41 = 6
    42 = 5
}

# whether an odd cycle was found ("yes"/"no")
$I = @{
    23 = "yes"
    24 = "yes"
    25 = "yes"
    26 = "yes"
    27 = "no"
    28 = "yes"
    29 = "yes"
    30 = "yes"
    31 = "yes"
    32 = "yes"
    33 = "yes"
    34 = "no"
    35 = "yes"
    36 = "yes"
    37 = "yes"
    38 = "yes"
    39 = "yes"
    40 = "yes"
    41 = "yes"
    42 = "yes"
}

# rows where the max clique exactly matches the goal (chromatic number) -
# these get the highlighted fill that already marks "goal" matches in column E
$Highlight = @(23, 27, 28, 32, 34)

for ($r = 23; $r -le 42; $r++) {
    $gAddr = "G" + $r
    $iAddr = "I" + $r
    $eAddr = "E" + $r

    $ws.Range($gAddr).Value = $G[$r]
    $ws.Range($iAddr).Value = $I[$r]

    if ($Highlight -contains $r) {
        # copy the highlighted fill/format from the goal column onto the
        # goal cell itself (upgrading it from the plain style), and onto
        # the new max-clique / odd-cycle cells for this row
        $ws.Range("E27").Copy()
        $ws.Range($eAddr).PasteSpecial(-4122)
        $ws.Range($gAddr).PasteSpecial(-4122)

        if ($I[$r] -eq "no") {
            $ws.Range($iAddr).PasteSpecial(-4122)
        }
    }
}

$excel.CutCopyMode = 0

# --- selection: mirror the new focus on the odd-cycle column ------------
$ws.Range("I35:I38").Select()
$ws.Range("I38").Activate()
